# Changed date and time fields for JGI app and verified database persistence
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# FOL_date (row 2): type changes from "date" to "text"
$ws.Range("C2").Value = "text"

# FOL_time_begin (row 5): type changes from "time" to "text"
$ws.Range("C5").Value = "text"

# FOL_time_end (row 6): type changes from "time" to "text"
$ws.Range("C6").Value = "text"

# Update the active selection to match the new cursor position
$ws.Activate()
$ws.Range("C7").Select()
